# Auto-generated edit script: Add data for 2022-12-09
# Updates 2022 full-year totals (column I) across Citywide Totals, By Neighborhood summary, and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6891
$ws.Range("I3").Value = 7149
$ws.Range("I4").Value = 1644
$ws.Range("I5").Value = 670
$ws.Range("I6").Value = 8360
$ws.Range("I7").Value = 24714

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 193
$ws.Range("I6").Value = 177
$ws.Range("I7").Value = 775
$ws.Range("I8").Value = 1479
$ws.Range("I9").Value = 128
$ws.Range("I11").Value = 377
$ws.Range("I14").Value = 137
$ws.Range("I15").Value = 286
$ws.Range("I19").Value = 696
$ws.Range("I20").Value = 608
$ws.Range("I21").Value = 110
$ws.Range("I23").Value = 245
$ws.Range("I26").Value = 34
$ws.Range("I27").Value = 215
$ws.Range("I33").Value = 1093
$ws.Range("I36").Value = 337
$ws.Range("I37").Value = 763
$ws.Range("I47").Value = 179
$ws.Range("I48").Value = 311
$ws.Range("I51").Value = 290
$ws.Range("I52").Value = 558
$ws.Range("I53").Value = 277
$ws.Range("I63").Value = 75
$ws.Range("I65").Value = 579
$ws.Range("I66").Value = 71
$ws.Range("I67").Value = 936
$ws.Range("I70").Value = 43
$ws.Range("I72").Value = 98
$ws.Range("I76").Value = 351
$ws.Range("I77").Value = 148
$ws.Range("I78").Value = 331
$ws.Range("I83").Value = 532
$ws.Range("I85").Value = 1107
$ws.Range("I86").Value = 158
$ws.Range("I88").Value = 226
$ws.Range("I89").Value = 292
$ws.Range("I90").Value = 319
$ws.Range("I95").Value = 377
$ws.Range("I97").Value = 221
$ws.Range("I99").Value = 433
$ws.Range("I101").Value = 24714

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 315
$ws.Range("I3").Value = 417
$ws.Range("I4").Value = 51
$ws.Range("I6").Value = 289
$ws.Range("I7").Value = 1107

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 182
$ws.Range("I7").Value = 558

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 149
$ws.Range("I4").Value = 37
$ws.Range("I7").Value = 377

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 436
$ws.Range("I3").Value = 427
$ws.Range("I4").Value = 93
$ws.Range("I7").Value = 1479

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I3").Value = 56
$ws.Range("I7").Value = 277

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 254
$ws.Range("I3").Value = 237
$ws.Range("I7").Value = 775

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 73
$ws.Range("I6").Value = 102
$ws.Range("I7").Value = 292

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 227
$ws.Range("I7").Value = 763

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 122
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 433

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 223
$ws.Range("I7").Value = 936

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 186
$ws.Range("I3").Value = 172
$ws.Range("I6").Value = 178
$ws.Range("I7").Value = 579

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I6").Value = 119
$ws.Range("I7").Value = 532

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 133
$ws.Range("I4").Value = 17
$ws.Range("I7").Value = 377

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 404
$ws.Range("I7").Value = 1093

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 224
$ws.Range("I6").Value = 222
$ws.Range("I7").Value = 696

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 49
$ws.Range("I6").Value = 161
$ws.Range("I7").Value = 311

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 159
$ws.Range("I7").Value = 351

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 67
$ws.Range("I7").Value = 177

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 85
$ws.Range("I7").Value = 331

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 88
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 83
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 172
$ws.Range("I6").Value = 210
$ws.Range("I7").Value = 608

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I5").Value = 12
$ws.Range("I7").Value = 337

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 56
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 84
$ws.Range("I7").Value = 286

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 34

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I2").Value = 21
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I2").Value = 44
$ws.Range("I7").Value = 128

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I3").Value = 60
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 221

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 43

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 77
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I2").Value = 60
$ws.Range("I7").Value = 215

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I4").Value = 76
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 319

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I4").Value = 29
$ws.Range("I7").Value = 290

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 98

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 148
